$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 37500
$ws.Range("J63").Value = 37500
$ws.Range("L63").Value = 37500
$ws.Range("N63").Value = -38748
$ws.Range("H66").Value = 37500
$ws.Range("J66").Value = 37500
$ws.Range("L66").Value = 112500
$ws.Range("N66").Value = -118740
$ws.Range("H96").Value = 1700
$ws.Range("I96").Value = 1700
$ws.Range("K96").Value = 5100
$ws.Range("M96").Value = -3727
$ws.Range("H98").Value = 4781.5894
$ws.Range("I98").Value = 3165.1
$ws.Range("J98").Value = 6646.769
$ws.Range("K98").Value = 3165.1
$ws.Range("L98").Value = 6646.769
$ws.Range("M98").Value = -1667.1
$ws.Range("N98").Value = -9642.769
$ws.Range("H100").Value = 25001152
$ws.Range("I100").Value = 28572602
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 28572602
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -28572061
$ws.Range("N100").Value = -2082
$ws.Range("H110").Value = 34950
$ws.Range("J110").Value = 34950
$ws.Range("L110").Value = 34950
$ws.Range("N110").Value = -43130
$ws.Range("H112").Value = 1273.1904
$ws.Range("J112").Value = 1273.1904
$ws.Range("L112").Value = 3819.5712
$ws.Range("N112").Value = -6035.5712
$ws.Range("H122").Value = 4781.5894
$ws.Range("I122").Value = 3165.1
$ws.Range("J122").Value = 6646.769
$ws.Range("K122").Value = 9495.299999999999
$ws.Range("L122").Value = 19940.307
$ws.Range("M122").Value = -7045.299999999999
$ws.Range("N122").Value = -24840.307
$ws.Range("H129").Value = 822.64
$ws.Range("I129").Value = 303.16666
$ws.Range("J129").Value = 855.79785
$ws.Range("K129").Value = 909.4999799999999
$ws.Range("L129").Value = 2567.39355
$ws.Range("M129").Value = 4090.50002
$ws.Range("N129").Value = -12567.39355
$ws.Range("H138").Value = 2808.3962
$ws.Range("I138").Value = 1940.8334
$ws.Range("J138").Value = 3062.3171
$ws.Range("K138").Value = 5822.5002
$ws.Range("L138").Value = 9186.951300000001
$ws.Range("M138").Value = -682.5002000000004
$ws.Range("N138").Value = -19466.9513

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 169.5
$ws.Range("I5").Value = 126
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 126
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -14
$ws.Range("N5").Value = -524
$ws.Range("H32").Value = 5569.698
$ws.Range("I32").Value = 5502.6387
$ws.Range("J32").Value = 5711.706
$ws.Range("K32").Value = 5502.6387
$ws.Range("L32").Value = 5711.706
$ws.Range("M32").Value = -5215.6387
$ws.Range("N32").Value = -6285.706
$ws.Range("H74").Value = 1706.1936
$ws.Range("I74").Value = 1141.1852
$ws.Range("J74").Value = 5520
$ws.Range("K74").Value = 1141.1852
$ws.Range("L74").Value = 5520
$ws.Range("M74").Value = -267.1851999999999
$ws.Range("N74").Value = -7268
$ws.Range("H77").Value = 1706.1936
$ws.Range("I77").Value = 1141.1852
$ws.Range("J77").Value = 5520
$ws.Range("K77").Value = 5705.925999999999
$ws.Range("L77").Value = 27600
$ws.Range("M77").Value = -1337.925999999999
$ws.Range("N77").Value = -36336
$ws.Range("H102").Value = 1871.2858
$ws.Range("I102").Value = 1849.8334
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1849.8334
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -227.8334
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 169.5
$ws.Range("I4").Value = 126
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 126
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -11
$ws.Range("N4").Value = -530
$ws.Range("H94").Value = 750
$ws.Range("I94").Value = 750
$ws.Range("K94").Value = 750
$ws.Range("M94").Value = -299
$ws.Range("H99").Value = 2159.6
$ws.Range("I99").Value = 1046.5294
$ws.Range("J99").Value = 4524.875
$ws.Range("K99").Value = 1046.5294
$ws.Range("L99").Value = 4524.875
$ws.Range("M99").Value = 451.4706000000001
$ws.Range("N99").Value = -7520.875
$ws.Range("H103").Value = 35052.59
$ws.Range("J103").Value = 35052.59
$ws.Range("L103").Value = 35052.59
$ws.Range("N103").Value = -37396.59
$ws.Range("H107").Value = 1314.2
$ws.Range("I107").Value = 449
$ws.Range("J107").Value = 3333
$ws.Range("K107").Value = 449
$ws.Range("L107").Value = 3333
$ws.Range("M107").Value = 1471
$ws.Range("N107").Value = -7173
$ws.Range("H134").Value = 2573.7837
$ws.Range("I134").Value = 1533.1482
$ws.Range("J134").Value = 5383.5
$ws.Range("K134").Value = 4599.444600000001
$ws.Range("L134").Value = 16150.5
$ws.Range("M134").Value = -2064.444600000001
$ws.Range("N134").Value = -21220.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 18188.4
$ws.Range("J9").Value = 18188.4
$ws.Range("L9").Value = 18188.4
$ws.Range("N9").Value = -18524.4
$ws.Range("H16").Value = 10102936
$ws.Range("I16").Value = 18520098
$ws.Range("J16").Value = 2339.8
$ws.Range("K16").Value = 18520098
$ws.Range("L16").Value = 2339.8
$ws.Range("M16").Value = -18519811
$ws.Range("N16").Value = -2913.8
$ws.Range("H31").Value = 5079.156
$ws.Range("J31").Value = 8115.8975
$ws.Range("L31").Value = 8115.8975
$ws.Range("N31").Value = -8705.897499999999
$ws.Range("H34").Value = 5079.156
$ws.Range("J34").Value = 8115.8975
$ws.Range("L34").Value = 8115.8975
$ws.Range("N34").Value = -8519.897499999999
$ws.Range("H113").Value = 10102936
$ws.Range("I113").Value = 18520098
$ws.Range("J113").Value = 2339.8
$ws.Range("K113").Value = 18520098
$ws.Range("L113").Value = 2339.8
$ws.Range("M113").Value = -18517928
$ws.Range("N113").Value = -6679.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5682393
$ws.Range("I113").Value = 618.0909
$ws.Range("J113").Value = 11364168
$ws.Range("K113").Value = 1854.2727
$ws.Range("L113").Value = 34092504
$ws.Range("M113").Value = 315.7273
$ws.Range("N113").Value = -34096844
$ws.Range("H122").Value = 2492.1553
$ws.Range("I122").Value = 986.1818
$ws.Range("J122").Value = 2844.617
$ws.Range("K122").Value = 8875.636199999999
$ws.Range("L122").Value = 25601.553
$ws.Range("M122").Value = -6425.636199999999
$ws.Range("N122").Value = -30501.553

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6145.93
$ws.Range("I70").Value = 5762.533
$ws.Range("J70").Value = 7030.6924
$ws.Range("K70").Value = 5762.533
$ws.Range("L70").Value = 7030.6924
$ws.Range("M70").Value = -5492.533
$ws.Range("N70").Value = -7570.6924
$ws.Range("H73").Value = 6145.93
$ws.Range("I73").Value = 5762.533
$ws.Range("J73").Value = 7030.6924
$ws.Range("K73").Value = 5762.533
$ws.Range("L73").Value = 7030.6924
$ws.Range("M73").Value = -4826.533
$ws.Range("N73").Value = -8902.6924
$ws.Range("H133").Value = 40911.668
$ws.Range("J133").Value = 40911.668
$ws.Range("L133").Value = 40911.668
$ws.Range("N133").Value = -51031.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 20792.75
$ws.Range("I88").Value = 6085.5
$ws.Range("J88").Value = 35500
$ws.Range("K88").Value = 6085.5
$ws.Range("L88").Value = 35500
$ws.Range("M88").Value = -5657.5
$ws.Range("N88").Value = -36356
$ws.Range("H91").Value = 20792.75
$ws.Range("I91").Value = 6085.5
$ws.Range("J91").Value = 35500
$ws.Range("K91").Value = 6085.5
$ws.Range("L91").Value = 35500
$ws.Range("M91").Value = -4603.5
$ws.Range("N91").Value = -38464
$ws.Range("H132").Value = 6307.5713
$ws.Range("I132").Value = 5188.5
$ws.Range("J132").Value = 7799.6665
$ws.Range("K132").Value = 15565.5
$ws.Range("L132").Value = 23398.9995
$ws.Range("M132").Value = -13035.5
$ws.Range("N132").Value = -28458.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 26700
$ws.Range("J64").Value = 26700
$ws.Range("L64").Value = 26700
$ws.Range("N64").Value = -27196
$ws.Range("H67").Value = 26700
$ws.Range("J67").Value = 26700
$ws.Range("L67").Value = 26700
$ws.Range("N67").Value = -28416
$ws.Range("H96").Value = 1144808.6
$ws.Range("I96").Value = 525924.4399999999
$ws.Range("J96").Value = 2382577
$ws.Range("K96").Value = 525924.4399999999
$ws.Range("L96").Value = 2382577
$ws.Range("M96").Value = -524551.4399999999
$ws.Range("N96").Value = -2385323
$ws.Range("H113").Value = 299
$ws.Range("I113").Value = 299
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 897
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1273
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 47623390
$ws.Range("J132").Value = 47623390
$ws.Range("L132").Value = 142870170
$ws.Range("N132").Value = -142875230

Write-Output "Applied all Chocobo_Profits price/profit updates"